$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 42.934631
$ws.Range("H2").Value = 128.803893
$ws.Range("I2").Value = 0.1048104500939054
$ws.Range("J2").Value = 0.1048104500939054
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 6591.71170699927
$ws.Range("R2").Value = 59325.40536299344
$ws.Range("S2").Value = 0.03324804401928892
$ws.Range("T2").Value = 0.03324804401928893

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 42.934631
$ws.Range("H3").Value = 128.803893
$ws.Range("I3").Value = 0.1048104500939054
$ws.Range("J3").Value = 0.1048104500939054
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 7247.352946903051
$ws.Range("R3").Value = 65226.17652212746
$ws.Range("S3").Value = 0.03655504374472222
$ws.Range("T3").Value = 0.03655504374472222

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 42.934631
$ws.Range("H4").Value = 128.803893
$ws.Range("I4").Value = 0.1048104500939054
$ws.Range("J4").Value = 0.1048104500939054
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 2923.432906987357
$ws.Range("R4").Value = 26310.89616288621
$ws.Range("S4").Value = 0.01474555173214649
$ws.Range("T4").Value = 0.0147455517321465

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 42.934631
$ws.Range("H5").Value = 128.803893
$ws.Range("I5").Value = 0.1048104500939054
$ws.Range("J5").Value = 0.1048104500939054
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 4017.078840628663
$ws.Range("R5").Value = 36153.70956565797
$ws.Range("S5").Value = 0.02026181059774777
$ws.Range("T5").Value = 0.02026181059774778

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 145.5961606666667
$ws.Range("H6").Value = 436.788482
$ws.Range("I6").Value = 0.3554240196315627
$ws.Range("J6").Value = 0.3554240196315627
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 22353.23547465945
$ws.Range("R6").Value = 201179.1192719351
$ws.Range("S6").Value = 0.1127478551960722
$ws.Range("T6").Value = 0.1127478551960723

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 145.5961606666667
$ws.Range("H7").Value = 436.788482
$ws.Range("I7").Value = 0.3554240196315627
$ws.Range("J7").Value = 0.3554240196315627
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 24576.58862994157
$ws.Range("R7").Value = 221189.2976694741
$ws.Range("S7").Value = 0.123962263055984
$ws.Range("T7").Value = 0.123962263055984

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 145.5961606666667
$ws.Range("H8").Value = 436.788482
$ws.Range("I8").Value = 0.3554240196315627
$ws.Range("J8").Value = 0.3554240196315627
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 9913.689655885284
$ws.Range("R8").Value = 89223.20690296755
$ws.Range("S8").Value = 0.05000382369915431
$ws.Range("T8").Value = 0.05000382369915433

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 145.5961606666667
$ws.Range("H9").Value = 436.788482
$ws.Range("I9").Value = 0.3554240196315627
$ws.Range("J9").Value = 0.3554240196315627
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 13622.3659705108
$ws.Range("R9").Value = 122601.2937345972
$ws.Range("S9").Value = 0.06871007768035213
$ws.Range("T9").Value = 0.06871007768035213

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 207.2564646666667
$ws.Range("H10").Value = 621.769394
$ws.Range("I10").Value = 0.5059468974261112
$ws.Range("J10").Value = 0.5059468974261113
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 31819.88135625405
$ws.Range("R10").Value = 286378.9322062865
$ws.Range("S10").Value = 0.1604968273867203
$ws.Range("T10").Value = 0.1604968273867204

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 207.2564646666667
$ws.Range("H11").Value = 621.769394
$ws.Range("I11").Value = 0.5059468974261112
$ws.Range("J11").Value = 0.5059468974261113
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 34984.82961147785
$ws.Range("R11").Value = 314863.4665033006
$ws.Range("S11").Value = 0.1764605623899848
$ws.Range("T11").Value = 0.1764605623899848

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 207.2564646666667
$ws.Range("H12").Value = 621.769394
$ws.Range("I12").Value = 0.5059468974261112
$ws.Range("J12").Value = 0.5059468974261113
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 14112.15969207691
$ws.Range("R12").Value = 127009.4372286922
$ws.Range("S12").Value = 0.07118055635703786
$ws.Range("T12").Value = 0.07118055635703789

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 207.2564646666667
$ws.Range("H13").Value = 621.769394
$ws.Range("I13").Value = 0.5059468974261112
$ws.Range("J13").Value = 0.5059468974261113
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 19391.46883074339
$ws.Range("R13").Value = 174523.2194766905
$ws.Range("S13").Value = 0.09780895129236825
$ws.Range("T13").Value = 0.09780895129236827

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.85349
$ws.Range("H14").Value = 41.56047
$ws.Range("I14").Value = 0.0338186328484206
$ws.Range("J14").Value = 0.03381863284842061
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 2126.91270633716
$ws.Range("R14").Value = 19142.21435703444
$ws.Range("S14").Value = 0.01072797027976737
$ws.Range("T14").Value = 0.01072797027976737

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.85349
$ws.Range("H15").Value = 41.56047
$ws.Range("I15").Value = 0.0338186328484206
$ws.Range("J15").Value = 0.03381863284842061
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 2338.464992895641
$ws.Range("R15").Value = 21046.18493606076
$ws.Range("S15").Value = 0.01179502236707407
$ws.Range("T15").Value = 0.01179502236707407

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.85349
$ws.Range("H16").Value = 41.56047
$ws.Range("I16").Value = 0.0338186328484206
$ws.Range("J16").Value = 0.03381863284842061
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 943.2886133951
$ws.Range("R16").Value = 8489.597520555901
$ws.Range("S16").Value = 0.004757869083951696
$ws.Range("T16").Value = 0.004757869083951698

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 13.85349
$ws.Range("H17").Value = 41.56047
$ws.Range("I17").Value = 0.0338186328484206
$ws.Range("J17").Value = 0.03381863284842061
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 1296.16955477877
$ws.Range("R17").Value = 11665.52599300893
$ws.Range("S17").Value = 0.006537771117627465
$ws.Range("T17").Value = 0.006537771117627466
